# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 34df19c).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1430
    3  = 7617
    5  = 325
    7  = 21
    9  = 5781
    12 = 21
    13 = 1768
    14 = 1277
    16 = 38
    17 = 7
    18 = 5515
    19 = 64
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
